$wb = $excel.ActiveWorkbook

# --- mmWave sheet: append rows 59-69 (new sensor events) ---
$ws = $wb.Worksheets.Item("mmWave")
$ws.Range("A59:A69").NumberFormat = "@"
$ws.Cells.Item(59, 1).Value = '2026-01-30'
$ws.Cells.Item(59, 2).Value = '15:09:23'
$ws.Cells.Item(59, 3).Value = '15:00'
$ws.Cells.Item(59, 4).Value = 'Living Room'
$ws.Cells.Item(59, 5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(59, 6).Value = 'Active'
$ws.Cells.Item(60, 1).Value = '2026-01-30'
$ws.Cells.Item(60, 2).Value = '15:09:33'
$ws.Cells.Item(60, 3).Value = '15:00'
$ws.Cells.Item(60, 4).Value = 'Living Room'
$ws.Cells.Item(60, 5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(60, 6).Value = 'Active'
$ws.Cells.Item(61, 1).Value = '2026-01-30'
$ws.Cells.Item(61, 2).Value = '15:09:44'
$ws.Cells.Item(61, 3).Value = '15:00'
$ws.Cells.Item(61, 4).Value = 'Living Room'
$ws.Cells.Item(61, 5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(61, 6).Value = 'Active'
$ws.Cells.Item(62, 1).Value = '2026-01-30'
$ws.Cells.Item(62, 2).Value = '15:09:54'
$ws.Cells.Item(62, 3).Value = '15:00'
$ws.Cells.Item(62, 4).Value = 'Living Room'
$ws.Cells.Item(62, 5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(62, 6).Value = 'Active'
$ws.Cells.Item(63, 1).Value = '2026-01-30'
$ws.Cells.Item(63, 2).Value = '15:10:22'
$ws.Cells.Item(63, 3).Value = '15:00'
$ws.Cells.Item(63, 4).Value = 'Living Room'
$ws.Cells.Item(63, 5).Value = 'FALL_DETECTED'
$ws.Cells.Item(63, 6).Value = 'EMERGENCY'
$ws.Cells.Item(64, 1).Value = '2026-01-30'
$ws.Cells.Item(64, 2).Value = '15:10:22'
$ws.Cells.Item(64, 3).Value = '15:00'
$ws.Cells.Item(64, 4).Value = 'Living Room'
$ws.Cells.Item(64, 5).Value = 'FALL_DETECTED'
$ws.Cells.Item(64, 6).Value = 'EMERGENCY'
$ws.Cells.Item(65, 1).Value = '2026-01-30'
$ws.Cells.Item(65, 2).Value = '15:10:38'
$ws.Cells.Item(65, 3).Value = '15:00'
$ws.Cells.Item(65, 4).Value = 'Living Room'
$ws.Cells.Item(65, 5).Value = 'CRITICAL EMERGENCY'
$ws.Cells.Item(65, 6).Value = 'FALL_DETECTED'
$ws.Cells.Item(66, 1).Value = '2026-01-30'
$ws.Cells.Item(66, 2).Value = '15:10:41'
$ws.Cells.Item(66, 3).Value = '15:00'
$ws.Cells.Item(66, 4).Value = 'Living Room'
$ws.Cells.Item(66, 5).Value = 'CRITICAL EMERGENCY'
$ws.Cells.Item(66, 6).Value = 'FALL_DETECTED'
$ws.Cells.Item(67, 1).Value = '2026-01-30'
$ws.Cells.Item(67, 2).Value = '15:10:43'
$ws.Cells.Item(67, 3).Value = '15:00'
$ws.Cells.Item(67, 4).Value = 'Living Room'
$ws.Cells.Item(67, 5).Value = 'FALL_DETECTED'
$ws.Cells.Item(67, 6).Value = 'EMERGENCY'
$ws.Cells.Item(68, 1).Value = '2026-01-30'
$ws.Cells.Item(68, 2).Value = '15:10:46'
$ws.Cells.Item(68, 3).Value = '15:00'
$ws.Cells.Item(68, 4).Value = 'Living Room'
$ws.Cells.Item(68, 5).Value = 'CRITICAL EMERGENCY'
$ws.Cells.Item(68, 6).Value = 'FALL_DETECTED'
$ws.Cells.Item(69, 1).Value = '2026-01-30'
$ws.Cells.Item(69, 2).Value = '15:11:23'
$ws.Cells.Item(69, 3).Value = '15:00'
$ws.Cells.Item(69, 4).Value = 'Living Room'
$ws.Cells.Item(69, 5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(69, 6).Value = 'Active'

# --- Proximity sheet: append rows 13-21 (new door events) ---
$ws = $wb.Worksheets.Item("Proximity")
$ws.Range("A13:A21").NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = '2026-01-30'
$ws.Cells.Item(13, 2).Value = '15:10:42'
$ws.Cells.Item(13, 3).Value = '15:00'
$ws.Cells.Item(13, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(13, 5).Value = 'ENTER'
$ws.Cells.Item(13, 6).Value = 'User ENTERED Living Room Main Door'
$ws.Cells.Item(14, 1).Value = '2026-01-30'
$ws.Cells.Item(14, 2).Value = '15:10:51'
$ws.Cells.Item(14, 3).Value = '15:00'
$ws.Cells.Item(14, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(14, 5).Value = 'EXIT'
$ws.Cells.Item(14, 6).Value = 'User EXITED Living Room Main Door'
$ws.Cells.Item(15, 1).Value = '2026-01-30'
$ws.Cells.Item(15, 2).Value = '15:11:00'
$ws.Cells.Item(15, 3).Value = '15:00'
$ws.Cells.Item(15, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(15, 5).Value = 'ENTER'
$ws.Cells.Item(15, 6).Value = 'User ENTERED Living Room Main Door'
$ws.Cells.Item(16, 1).Value = '2026-01-30'
$ws.Cells.Item(16, 2).Value = '15:11:03'
$ws.Cells.Item(16, 3).Value = '15:00'
$ws.Cells.Item(16, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(16, 5).Value = 'EXIT'
$ws.Cells.Item(16, 6).Value = 'User EXITED Living Room Main Door'
$ws.Cells.Item(17, 1).Value = '2026-01-30'
$ws.Cells.Item(17, 2).Value = '15:11:06'
$ws.Cells.Item(17, 3).Value = '15:00'
$ws.Cells.Item(17, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(17, 5).Value = 'ENTER'
$ws.Cells.Item(17, 6).Value = 'User ENTERED Living Room Main Door'
$ws.Cells.Item(18, 1).Value = '2026-01-30'
$ws.Cells.Item(18, 2).Value = '15:11:10'
$ws.Cells.Item(18, 3).Value = '15:00'
$ws.Cells.Item(18, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(18, 5).Value = 'EXIT'
$ws.Cells.Item(18, 6).Value = 'User EXITED Living Room Main Door'
$ws.Cells.Item(19, 1).Value = '2026-01-30'
$ws.Cells.Item(19, 2).Value = '15:11:14'
$ws.Cells.Item(19, 3).Value = '15:00'
$ws.Cells.Item(19, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(19, 5).Value = 'ENTER'
$ws.Cells.Item(19, 6).Value = 'User ENTERED Living Room Main Door'
$ws.Cells.Item(20, 1).Value = '2026-01-30'
$ws.Cells.Item(20, 2).Value = '15:11:18'
$ws.Cells.Item(20, 3).Value = '15:00'
$ws.Cells.Item(20, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(20, 5).Value = 'EXIT'
$ws.Cells.Item(20, 6).Value = 'User EXITED Living Room Main Door'
$ws.Cells.Item(21, 1).Value = '2026-01-30'
$ws.Cells.Item(21, 2).Value = '15:11:21'
$ws.Cells.Item(21, 3).Value = '15:00'
$ws.Cells.Item(21, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(21, 5).Value = 'ENTER'
$ws.Cells.Item(21, 6).Value = 'User ENTERED Living Room Main Door'

# --- Camera sheet: append rows 16-22 (new image-capture events) ---
$ws = $wb.Worksheets.Item("Camera")
$ws.Range("A16:A22").NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = '2026-01-30'
$ws.Cells.Item(16, 2).Value = '15:10:43'
$ws.Cells.Item(16, 3).Value = '15:00'
$ws.Cells.Item(16, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(16, 5).Value = 'Image Captured (ENTER)'
$ws.Cells.Item(16, 6).Value = 'Active'
$ws.Cells.Item(17, 1).Value = '2026-01-30'
$ws.Cells.Item(17, 2).Value = '15:10:50'
$ws.Cells.Item(17, 3).Value = '15:00'
$ws.Cells.Item(17, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(17, 5).Value = 'Image Captured (ENTER)'
$ws.Cells.Item(17, 6).Value = 'Active'
$ws.Cells.Item(18, 1).Value = '2026-01-30'
$ws.Cells.Item(18, 2).Value = '15:11:01'
$ws.Cells.Item(18, 3).Value = '15:00'
$ws.Cells.Item(18, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(18, 5).Value = 'Image Captured (ENTER)'
$ws.Cells.Item(18, 6).Value = 'Active'
$ws.Cells.Item(19, 1).Value = '2026-01-30'
$ws.Cells.Item(19, 2).Value = '15:11:07'
$ws.Cells.Item(19, 3).Value = '15:00'
$ws.Cells.Item(19, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(19, 5).Value = 'Image Captured (ENTER)'
$ws.Cells.Item(19, 6).Value = 'Active'
$ws.Cells.Item(20, 1).Value = '2026-01-30'
$ws.Cells.Item(20, 2).Value = '15:11:14'
$ws.Cells.Item(20, 3).Value = '15:00'
$ws.Cells.Item(20, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(20, 5).Value = 'Image Captured (ENTER)'
$ws.Cells.Item(20, 6).Value = 'Active'
$ws.Cells.Item(21, 1).Value = '2026-01-30'
$ws.Cells.Item(21, 2).Value = '15:11:18'
$ws.Cells.Item(21, 3).Value = '15:00'
$ws.Cells.Item(21, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(21, 5).Value = 'Image Captured (ENTER)'
$ws.Cells.Item(21, 6).Value = 'Active'
$ws.Cells.Item(22, 1).Value = '2026-01-30'
$ws.Cells.Item(22, 2).Value = '15:11:22'
$ws.Cells.Item(22, 3).Value = '15:00'
$ws.Cells.Item(22, 4).Value = 'Living Room Main Door'
$ws.Cells.Item(22, 5).Value = 'Image Captured (ENTER)'
$ws.Cells.Item(22, 6).Value = 'Active'
